$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as the new row 33; every
# existing data row from the old row 33 onward shifts down by one
# (old row 33 -> new row 34, ..., old row 95 -> new row 96).
$ws.Rows.Item(33).Insert()

$ws.Range("A33").Value = 10
$ws.Range("B33").Value = "Vega Modelo de Temuco"
$ws.Range("C33").Value = "La Araucanía"
$ws.Range("D33").Value = 44526
$ws.Range("E33").Value = 9
$ws.Range("F33").Value = 100112031
$ws.Range("G33").Value = "Poroto verde"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 20
$ws.Range("K33").Value = 35000
$ws.Range("L33").Value = 35000
$ws.Range("M33").Value = 35000
$ws.Range("N33").Value = '$/malla 25 kilos'
$ws.Range("O33").Value = "Provincia de Limarí"
$ws.Range("P33").Value = 1400
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"
